$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B7").Value = "417823 - PREMIUM SAÚDE S.A"

$ws.Range("C2").Value = "13/04/2023  10:06:53"
$ws.Range("C3").Value = "13/04/2023  15:44:44"
$ws.Range("C4").Value = "13/04/2023  17:09:35"
$ws.Range("C5").Value = "14/04/2023  10:46:04"
$ws.Range("C6").Value = "14/04/2023  11:24:10"
$ws.Range("C7").Value = "14/04/2023  15:09:48"

$ws.Range("D2").Value = 12168346
$ws.Range("E2").Value = 8606403
$ws.Range("D3").Value = 12169638
$ws.Range("E3").Value = 8608000
$ws.Range("D4").Value = 12169885
$ws.Range("E4").Value = 8608362
$ws.Range("D5").Value = 12170717
$ws.Range("E5").Value = 8609371
$ws.Range("D6").Value = 12170838
$ws.Range("E6").Value = 8609527
$ws.Range("D7").Value = 12171665
$ws.Range("E7").Value = 8610551

$ws.Range("F2").Value = "MELINDA GOULART CRUZ"
$ws.Range("F3").Value = "JORGE EUSTACIO MEDEIROS"
$ws.Range("F4").Value = "NATALIA SOLANO ROSSELIS PEREIRA DA SILVA"
$ws.Range("F5").Value = "RAFAELLE FONSECA DE OLIVEIRA"
$ws.Range("F6").Value = "VANESSA WINTER NUNES FORTES"
$ws.Range("F7").Value = "OLIVER MIRANDA PORTO"

$ws.Range("I4").Value = "No Assistencial"

$ws.Range("A8:J9").Delete()
